# Apply odds updates to the Jogos da Semana FlashScore worksheet
# per the commit "Atualizando o arquivo XLSX".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.1
$ws.Range("I2").Value = 3.2
$ws.Range("L2").Value = 3.6
$ws.Range("S2").Value = 2.02
$ws.Range("T2").Value = 1.88
$ws.Range("Y2").Value = 1.53
$ws.Range("Z2").Value = 2.38
$ws.Range("AC2").Value = 10
$ws.Range("G3").Value = 2.35
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 2.88
$ws.Range("O3").Value = 1.22
$ws.Range("P3").Value = 4.33
$ws.Range("Y3").Value = 1.57
$ws.Range("Z3").Value = 2.25
$ws.Range("AH3").Value = 7
$ws.Range("AL3").Value = 13
$ws.Range("AM3").Value = 19
$ws.Range("J4").Value = 3.4
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.73
$ws.Range("AG4").Value = 9
$ws.Range("AH4").Value = 6.5
$ws.Range("AK4").Value = 351
$ws.Range("AL4").Value = 9
$ws.Range("H6").Value = 3.7
$ws.Range("I6").Value = 3.35
$ws.Range("J6").Value = 2.55
$ws.Range("K6").Value = 2.2
$ws.Range("L6").Value = 3.8
$ws.Range("O6").Value = 1.26
$ws.Range("P6").Value = 3.5
$ws.Range("Q6").Value = 1.78
$ws.Range("R6").Value = 1.93
$ws.Range("U6").Value = 2.87
$ws.Range("V6").Value = 1.37
$ws.Range("X6").Value = 2.85
$ws.Range("Y6").Value = 1.7
$ws.Range("Z6").Value = 2.02
$ws.Range("AA6").Value = 8
$ws.Range("AB6").Value = 9.75
$ws.Range("AC6").Value = 8.5
$ws.Range("AD6").Value = 17
$ws.Range("AE6").Value = 15
$ws.Range("AF6").Value = 25
$ws.Range("AH6").Value = 7
$ws.Range("AI6").Value = 14
$ws.Range("AJ6").Value = 60
$ws.Range("AK6").Value = 450
$ws.Range("AL6").Value = 11
$ws.Range("AM6").Value = 18
$ws.Range("AO6").Value = 45
$ws.Range("AP6").Value = 28
$ws.Range("AQ6").Value = 35
$ws.Range("G7").Value = 1.57
$ws.Range("L7").Value = 7.5
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("O7").Value = 1.44
$ws.Range("P7").Value = 2.63
$ws.Range("Y7").Value = 2.63
$ws.Range("Z7").Value = 1.44
$ws.Range("AH7").Value = 7.5
$ws.Range("AP7").Value = 67
$ws.Range("AR7").Value = 1.88
$ws.Range("AS7").Value = 1.98
$ws.Range("G8").Value = 2.8
$ws.Range("I8").Value = 2.63
$ws.Range("L8").Value = 3.25
$ws.Range("AM8").Value = 12
$ws.Range("Q9").Value = 1.8
$ws.Range("R9").Value = 2
$ws.Range("Y9").Value = 1.62
$ws.Range("Z9").Value = 2.2
$ws.Range("AA9").Value = 11
$ws.Range("AG9").Value = 12
$ws.Range("AH9").Value = 7
$ws.Range("AK9").Value = 151
$ws.Range("AL9").Value = 9
$ws.Range("AQ9").Value = 23
$ws.Range("L10").Value = 6.5
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 8
$ws.Range("Y10").Value = 2.2
$ws.Range("Z10").Value = 1.62
$ws.Range("AG10").Value = 8
$ws.Range("AH10").Value = 7
$ws.Range("AJ10").Value = 81
$ws.Range("AL10").Value = 12
$ws.Range("Q11").Value = 1.98
$ws.Range("R11").Value = 1.88
$ws.Range("G14").Value = 1.42
$ws.Range("H14").Value = 4.5
$ws.Range("M14").Value = 1.04
$ws.Range("N14").Value = 13
$ws.Range("O14").Value = 1.2
$ws.Range("P14").Value = 4.33
$ws.Range("Q14").Value = 1.65
$ws.Range("R14").Value = 2.2
$ws.Range("S14").Value = 2.03
$ws.Range("T14").Value = 1.83
$ws.Range("U14").Value = 2.63
$ws.Range("V14").Value = 1.44
$ws.Range("Y14").Value = 1.91
$ws.Range("Z14").Value = 1.91
$ws.Range("AA14").Value = 7.5
$ws.Range("AB14").Value = 7
$ws.Range("AG14").Value = 13
$ws.Range("AH14").Value = 8.5
$ws.Range("AI14").Value = 17
$ws.Range("AL14").Value = 21
$ws.Range("AQ14").Value = 51
$ws.Range("I16").Value = 11
$ws.Range("N16").Value = 17
$ws.Range("AB16").Value = 9
$ws.Range("AC16").Value = 12
$ws.Range("AL16").Value = 41
$ws.Range("AN16").Value = 34
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 5.75
$ws.Range("I17").Value = 1.33
$ws.Range("J17").Value = 6
$ws.Range("L17").Value = 1.73
$ws.Range("M17").Value = 1.02
$ws.Range("N17").Value = 12
$ws.Range("O17").Value = 1.13
$ws.Range("P17").Value = 5.5
$ws.Range("Q17").Value = 1.4
$ws.Range("R17").Value = 2.75
$ws.Range("U17").Value = 2
$ws.Range("V17").Value = 1.73
$ws.Range("Y17").Value = 1.73
$ws.Range("Z17").Value = 2
$ws.Range("AB17").Value = 41
$ws.Range("AC17").Value = 19
$ws.Range("AD17").Value = 67
$ws.Range("AF17").Value = 41
$ws.Range("AG17").Value = 23
$ws.Range("AI17").Value = 21
$ws.Range("AJ17").Value = 51
$ws.Range("AK17").Value = 151
$ws.Range("AL17").Value = 10
$ws.Range("AM17").Value = 8.5
$ws.Range("AO17").Value = 9.5
$ws.Range("AQ17").Value = 21
$ws.Range("G18").Value = 2.3
$ws.Range("H18").Value = 3.4
$ws.Range("I18").Value = 2.8
$ws.Range("J18").Value = 2.88
$ws.Range("L18").Value = 3.4
$ws.Range("M18").Value = 1.02
$ws.Range("N18").Value = 12
$ws.Range("O18").Value = 1.22
$ws.Range("P18").Value = 4
$ws.Range("Q18").Value = 1.8
$ws.Range("R18").Value = 2
$ws.Range("U18").Value = 2.75
$ws.Range("V18").Value = 1.4
$ws.Range("Y18").Value = 1.67
$ws.Range("Z18").Value = 2.1
$ws.Range("AA18").Value = 9.5
$ws.Range("AB18").Value = 12
$ws.Range("AC18").Value = 9.5
$ws.Range("AD18").Value = 21
$ws.Range("AG18").Value = 12
$ws.Range("AH18").Value = 6.5
$ws.Range("AI18").Value = 13
$ws.Range("AM18").Value = 15
$ws.Range("AN18").Value = 11
$ws.Range("AO18").Value = 29
$ws.Range("AP18").Value = 21
$ws.Range("AQ18").Value = 29
$ws.Range("H19").Value = 9
$ws.Range("J19").Value = 1.44
$ws.Range("K19").Value = 3.25
$ws.Range("N19").Value = 23
$ws.Range("U19").Value = 1.91
$ws.Range("V19").Value = 1.91
$ws.Range("AA19").Value = 11
$ws.Range("AL19").Value = 41
$ws.Range("AO19").Value = 201
$ws.Range("AQ19").Value = 67
$ws.Range("M20").Value = 1.07
$ws.Range("N20").Value = 7.5
$ws.Range("G21").Value = 3.1
$ws.Range("I21").Value = 2.38
$ws.Range("M21").Value = 1.07
$ws.Range("N21").Value = 7.5
$ws.Range("Q21").Value = 2.08
$ws.Range("R21").Value = 1.73
$ws.Range("AO21").Value = 21
$ws.Range("Q22").Value = 2
$ws.Range("R22").Value = 1.85
$ws.Range("U22").Value = 3.4
$ws.Range("V22").Value = 1.3
$ws.Range("G23").Value = 4.3
$ws.Range("H23").Value = 3.6
$ws.Range("I23").Value = 1.78
$ws.Range("J23").Value = 4.45
$ws.Range("L23").Value = 2.37
$ws.Range("R23").Value = 1.93
$ws.Range("U23").Value = 2.92
$ws.Range("V23").Value = 1.37
$ws.Range("X23").Value = 2.82
$ws.Range("Y23").Value = 1.72
$ws.Range("Z23").Value = 2
$ws.Range("AA23").Value = 12.5
$ws.Range("AB23").Value = 28
$ws.Range("AC23").Value = 14.5
$ws.Range("AD23").Value = 80
$ws.Range("AE23").Value = 40
$ws.Range("AF23").Value = 45
$ws.Range("AH23").Value = 7.3
$ws.Range("AI23").Value = 15
$ws.Range("AK23").Value = 500
$ws.Range("AL23").Value = 7.1
$ws.Range("AM23").Value = 9.25
$ws.Range("AN23").Value = 8.5
$ws.Range("AO23").Value = 15.5
$ws.Range("AP23").Value = 15
$ws.Range("AQ23").Value = 27
